$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared-string insertion order matters (matches the original authoring
# order: row labels for "Unclosed", "BigIntegers" and "BigDecimals" were
# entered before the BigInteger/BigDecimal expression values).
$ws.Range("A19").Value = "Unclosed:"
$ws.Range("B19").Value = 'Unclosed "double quote'

$ws.Range("A20").Value = "BigIntegers:"
$ws.Range("A24").Value = "BigDecimals:"

$ws.Range("B20").Value = '${biZero}'
$ws.Range("B21").Value = '${biAnswer}'
$ws.Range("B22").Value = '${biBiggerThanLong}'
$ws.Range("B23").Value = '${biBiggerThanDouble}'

$ws.Range("B24").Value = '${bdZero}'
$ws.Range("B25").Value = '${bdAnswer}'
$ws.Range("B26").Value = '${bdSmallerThanNormal}'
$ws.Range("B27").Value = '${bdBiggerThanDouble}'

$ws.Range("A28").Value = "ValueHolder:"
$ws.Range("B28").Value = '${valueHolder.answer}'
$ws.Range("B29").Value = '${valueHolder.IHaveAQuestion}'

$ws.Columns.Item(1).ColumnWidth = 11.6
